# "4.0.3 model and data"
#
# The "Boolean" input-data-type sheet listed single combined CSV files for
# two transportation categories ("trans/BVTQaZ/BVTQaZ.csv" and
# "trans/VTQaZ/VTQaZ.csv"). Each of those is split into six per-vehicle-type
# CSV files (LDVs, HDVs, aircraft, rail, ships, motorbikes), so each of the
# two rows becomes six rows in place.

$wb = $excel.ActiveWorkbook

$wsBool = $wb.Worksheets.Item("Boolean")

# Row 17 currently holds "trans/BVTQaZ/BVTQaZ.csv" -> expand into 6 rows.
$wsBool.Rows.Item(18).Resize(5).Insert()
$wsBool.Range("A17").Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBool.Range("A18").Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBool.Range("A19").Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBool.Range("A20").Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBool.Range("A21").Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBool.Range("A22").Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# Original row 21 "trans/VTQaZ/VTQaZ.csv" has shifted down to row 26 because
# of the 5 rows just inserted above it -> expand into 6 rows too.
$wsBool.Rows.Item(27).Resize(5).Insert()
$wsBool.Range("A26").Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBool.Range("A27").Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBool.Range("A28").Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBool.Range("A29").Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBool.Range("A30").Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBool.Range("A31").Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# A handful of extra blank rows below the data picked up formatting in the
# source edit (selection dragged a little past the last row) - mirror that
# so the sheet's used range/dimension lines up.
$wsBool.Range("A33:A38").EntireRow.Font.Name = "Calibri"

# Leave the Boolean sheet scrolled to its last edited row.
$wsBool.Activate()
$wsBool.Range("A32").Select()

# The Integer sheet's selection moved to A13.
$wsInt = $wb.Worksheets.Item("Integer")
$wsInt.Activate()
$wsInt.Range("A13").Select()

# The About sheet ends up as the active tab.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
